$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "prix_tvac" (D column) values for the existing rows
$ws.Range("D2").Value = 11.35
$ws.Range("D3").Value = 14.21
$ws.Range("D4").Value = 11.35
$ws.Range("D5").Value = 11.35
$ws.Range("D6").Value = 7.13
$ws.Range("D7").Value = 9.0400000000000009
$ws.Range("D8").Value = 7.34
$ws.Range("D9").Value = 9.15

# Rows 10 and 11 previously had no value/format in column D; give them
# the same number format/style as the rest of the column before setting
# their new prices.
$ws.Range("D9").Copy()
$ws.Range("D10:D11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D10").Value = 16.64
$ws.Range("D11").Value = 21.91

# Reflect the selection that was active when the workbook was saved
$ws.Range("F2:F11").Select()
